$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H61").Value = 55
$ws.Range("I61").Value = 55
$ws.Range("K61").Value = 165
$ws.Range("M61").Value = 7
$ws.Range("H115").Value = 1185
$ws.Range("I115").Value = 1185
$ws.Range("K115").Value = 3555
$ws.Range("M115").Value = -1988

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7463.9165
$ws.Range("I61").Value = 3758.375
$ws.Range("J61").Value = 14875
$ws.Range("K61").Value = 3758.375
$ws.Range("L61").Value = 14875
$ws.Range("M61").Value = -3546.375
$ws.Range("N61").Value = -15299
$ws.Range("H74").Value = 9000
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126
$ws.Range("H77").Value = 9000
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632
$ws.Range("H97").Value = 2245.3333
$ws.Range("I97").Value = 2276
$ws.Range("K97").Value = 2276
$ws.Range("M97").Value = -1780
$ws.Range("H122").Value = 2000
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 5047.2666
$ws.Range("I132").Value = 2428.0908
$ws.Range("J132").Value = 12250
$ws.Range("K132").Value = 7284.2724
$ws.Range("L132").Value = 36750
$ws.Range("M132").Value = -4754.2724
$ws.Range("N132").Value = -41810
$ws.Range("H136").Value = 7463.9165
$ws.Range("I136").Value = 3758.375
$ws.Range("J136").Value = 14875
$ws.Range("K136").Value = 11275.125
$ws.Range("L136").Value = 44625
$ws.Range("M136").Value = -8725.125
$ws.Range("N136").Value = -49725

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4000
$ws.Range("I94").Value = 4000
$ws.Range("K94").Value = 4000
$ws.Range("M94").Value = -3549
$ws.Range("H99").Value = 2339.8
$ws.Range("I99").Value = 2339.8
$ws.Range("K99").Value = 2339.8
$ws.Range("M99").Value = -841.8000000000002
$ws.Range("H134").Value = 5078.6924
$ws.Range("I134").Value = 2365.7273
$ws.Range("K134").Value = 7097.1819
$ws.Range("M134").Value = -4562.1819

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4238.5
$ws.Range("J62").Value = 3999.6667
$ws.Range("L62").Value = 3999.6667
$ws.Range("N62").Value = -5247.6667
$ws.Range("H65").Value = 4238.5
$ws.Range("J65").Value = 3999.6667
$ws.Range("L65").Value = 19998.3335
$ws.Range("N65").Value = -26238.3335
$ws.Range("H94").Value = 1999.5
$ws.Range("J94").Value = 1999.5
$ws.Range("L94").Value = 1999.5
$ws.Range("N94").Value = -2901.5
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H132").Value = 2719.7058
$ws.Range("I132").Value = 1710.6154
$ws.Range("K132").Value = 5131.8462
$ws.Range("M132").Value = -2601.8462
$ws.Range("H134").Value = 7270.643
$ws.Range("I134").Value = 1865.7778
$ws.Range("K134").Value = 5597.3334
$ws.Range("M134").Value = -3062.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 324.22726
$ws.Range("J2").Value = 289.75
$ws.Range("L2").Value = 1738.5
$ws.Range("N2").Value = -1964.5
$ws.Range("H11").Value = 900
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 900
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2700
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -2980
$ws.Range("H12").Value = 30.363636
$ws.Range("I12").Value = 39
$ws.Range("J12").Value = 27.125
$ws.Range("K12").Value = 117
$ws.Range("L12").Value = 81.375
$ws.Range("M12").Value = 56
$ws.Range("N12").Value = -427.375
$ws.Range("H15").Value = 400.33334
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1199
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 3597
$ws.Range("M15").Value = 137
$ws.Range("N15").Value = -3877
$ws.Range("H26").Value = 220
$ws.Range("J26").Value = 350
$ws.Range("L26").Value = 1050
$ws.Range("N26").Value = -1626
$ws.Range("H92").Value = 849.75
$ws.Range("I92").Value = 866.6667
$ws.Range("K92").Value = 2600.0001
$ws.Range("M92").Value = -1352.0001
$ws.Range("H117").Value = 2579.4
$ws.Range("J117").Value = 2299.3333
$ws.Range("L117").Value = 6897.999899999999
$ws.Range("N117").Value = -13781.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H25").Value = 6662.6665
$ws.Range("J25").Value = 6662.6665
$ws.Range("L25").Value = 6662.6665
$ws.Range("N25").Value = -7720.6665
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H80").Value = 3240.889
$ws.Range("H83").Value = 3240.889

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 45062.5
$ws.Range("J6").Value = 45062.5
$ws.Range("L6").Value = 45062.5
$ws.Range("N6").Value = -45286.5
$ws.Range("H16").Value = 499
$ws.Range("I16").Value = 499
$ws.Range("K16").Value = 499
$ws.Range("M16").Value = -329
$ws.Range("H82").Value = 1670.8334
$ws.Range("I82").Value = 896.6667
$ws.Range("J82").Value = 2445
$ws.Range("K82").Value = 896.6667
$ws.Range("L82").Value = 2445
$ws.Range("M82").Value = -535.6667
$ws.Range("N82").Value = -3167
$ws.Range("H85").Value = 1670.8334
$ws.Range("I85").Value = 896.6667
$ws.Range("J85").Value = 2445
$ws.Range("K85").Value = 896.6667
$ws.Range("L85").Value = 2445
$ws.Range("M85").Value = 351.3333
$ws.Range("N85").Value = -4941
$ws.Range("H101").Value = 13615.6
$ws.Range("J101").Value = 13615.6
$ws.Range("L101").Value = 13615.6
$ws.Range("N101").Value = -20105.6
$ws.Range("H122").Value = 3005.3
$ws.Range("I122").Value = 2908.6667
$ws.Range("K122").Value = 8726.000100000001
$ws.Range("M122").Value = -6276.000100000001
$ws.Range("H132").Value = 14563
$ws.Range("I132").Value = 12100.8
$ws.Range("K132").Value = 36302.39999999999
$ws.Range("M132").Value = -33772.39999999999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 17199.7
$ws.Range("I136").Value = 13400
$ws.Range("K136").Value = 40200
$ws.Range("M136").Value = -37650

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4500
$ws.Range("J62").Value = 4500
$ws.Range("L62").Value = 4500
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 4500
$ws.Range("J65").Value = 4500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -28740
$ws.Range("H81").Value = 6500
$ws.Range("I81").Value = 2500
$ws.Range("J81").Value = 7642.857
$ws.Range("K81").Value = 5000
$ws.Range("L81").Value = 15285.714
$ws.Range("M81").Value = -3939
$ws.Range("N81").Value = -17407.714
$ws.Range("H84").Value = 6500
$ws.Range("I84").Value = 2500
$ws.Range("J84").Value = 7642.857
$ws.Range("K84").Value = 25000
$ws.Range("L84").Value = 76428.57000000001
$ws.Range("M84").Value = -19696
$ws.Range("N84").Value = -87036.57000000001
$ws.Range("H96").Value = 792
$ws.Range("I96").Value = 792
$ws.Range("K96").Value = 792
$ws.Range("M96").Value = 581
$ws.Range("H132").Value = 7663.2856
$ws.Range("I132").Value = 4229.375
$ws.Range("J132").Value = 12241.833
$ws.Range("K132").Value = 12688.125
$ws.Range("L132").Value = 36725.499
$ws.Range("M132").Value = -10158.125
$ws.Range("N132").Value = -41785.499
